{"js": "// Replace each \"before\" division expression with its \"after\" value.\n// Every cell text is unique in the document, so a plain body-wide search\n// for the exact old string unambiguously identifies the run to update.\nconst replacements = [\n  [\"77\u00f78=9, 5\", \"97\u00f72=48, 1\"],\n  [\"29\u00f73=9, 2\", \"54\u00f76=9, 0\"],\n  [\"51\u00f73=17, 0\", \"67\u00f73=22, 1\"],\n  [\"91\u00f72=45, 1\", \"96\u00f74=24, 0\"],\n  [\"75\u00f76=12, 3\", \"68\u00f75=13, 3\"],\n  [\"29\u00f76=4, 5\", \"21\u00f72=10, 1\"],\n  [\"24\u00f79=2, 6\", \"11\u00f73=3, 2\"],\n  [\"10\u00f72=5, 0\", \"28\u00f76=4, 4\"],\n  [\"85\u00f76=14, 1\", \"44\u00f79=4, 8\"],\n  [\"47\u00f73=15, 2\", \"14\u00f73=4, 2\"],\n  [\"18\u00f79=2, 0\", \"26\u00f76=4, 2\"],\n  [\"16\u00f75=3, 1\", \"25\u00f73=8, 1\"],\n  [\"67\u00f79=7, 4\", \"62\u00f75=12, 2\"],\n  [\"68\u00f76=11, 2\", \"49\u00f73=16, 1\"],\n  [\"12\u00f76=2, 0\", \"32\u00f73=10, 2\"],\n  [\"29\u00f77=4, 1\", \"78\u00f78=9, 6\"],\n  [\"93\u00f79=10, 3\", \"73\u00f76=12, 1\"],\n  [\"90\u00f73=30, 0\", \"36\u00f79=4, 0\"],\n  [\"79\u00f74=19, 3\", \"14\u00f77=2, 0\"],\n  [\"59\u00f75=11, 4\", \"57\u00f74=14, 1\"],\n  [\"52\u00f79=5, 7\", \"50\u00f73=16, 2\"],\n  [\"43\u00f77=6, 1\", \"98\u00f78=12, 2\"],\n  [\"12\u00f77=1, 5\", \"14\u00f72=7, 0\"],\n  [\"82\u00f72=41, 0\", \"97\u00f77=13, 6\"],\n  [\"13\u00f79=1, 4\", \"20\u00f77=2, 6\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"before\" division expression with its \"after\" value.\n# Every cell text is unique in the document, so Find/Execute with\n# wdReplaceAll unambiguously replaces exactly the one matching run.\n$replacements = @(\n  @(\"77\u00f78=9, 5\", \"97\u00f72=48, 1\"),\n  @(\"29\u00f73=9, 2\", \"54\u00f76=9, 0\"),\n  @(\"51\u00f73=17, 0\", \"67\u00f73=22, 1\"),\n  @(\"91\u00f72=45, 1\", \"96\u00f74=24, 0\"),\n  @(\"75\u00f76=12, 3\", \"68\u00f75=13, 3\"),\n  @(\"29\u00f76=4, 5\", \"21\u00f72=10, 1\"),\n  @(\"24\u00f79=2, 6\", \"11\u00f73=3, 2\"),\n  @(\"10\u00f72=5, 0\", \"28\u00f76=4, 4\"),\n  @(\"85\u00f76=14, 1\", \"44\u00f79=4, 8\"),\n  @(\"47\u00f73=15, 2\", \"14\u00f73=4, 2\"),\n  @(\"18\u00f79=2, 0\", \"26\u00f76=4, 2\"),\n  @(\"16\u00f75=3, 1\", \"25\u00f73=8, 1\"),\n  @(\"67\u00f79=7, 4\", \"62\u00f75=12, 2\"),\n  @(\"68\u00f76=11, 2\", \"49\u00f73=16, 1\"),\n  @(\"12\u00f76=2, 0\", \"32\u00f73=10, 2\"),\n  @(\"29\u00f77=4, 1\", \"78\u00f78=9, 6\"),\n  @(\"93\u00f79=10, 3\", \"73\u00f76=12, 1\"),\n  @(\"90\u00f73=30, 0\", \"36\u00f79=4, 0\"),\n  @(\"79\u00f74=19, 3\", \"14\u00f77=2, 0\"),\n  @(\"59\u00f75=11, 4\", \"57\u00f74=14, 1\"),\n  @(\"52\u00f79=5, 7\", \"50\u00f73=16, 2\"),\n  @(\"43\u00f77=6, 1\", \"98\u00f78=12, 2\"),\n  @(\"12\u00f77=1, 5\", \"14\u00f72=7, 0\"),\n  @(\"82\u00f72=41, 0\", \"97\u00f77=13, 6\"),\n  @(\"13\u00f79=1, 4\", \"20\u00f77=2, 6\")\n)\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindStop = 0\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = $wdFindStop\n\n  $null = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceAll)\n}\n"}
